$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price strings so Excel
# doesn't silently reinterpret them as numbers (matches source formatting).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"

# Apply updated cryptos list values
# Row 2
$ws.Range("D2").Value = '65.106.84'
$ws.Range("E2").Value = '  +0.12%  '

# Row 3
$ws.Range("D3").Value = '3.522.89'
$ws.Range("E3").Value = '  -1.33%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '593.24'
$ws.Range("E5").Value = '  -1.49%  '

# Row 6
$ws.Range("D6").Value = '134.48'
$ws.Range("E6").Value = '  -0.77%  '

# Row 7
$ws.Range("D7").Value = '3.524.80'

# Row 9
$ws.Range("E9").Value = '  -0.57%  '

# Row 10
$ws.Range("E10").Value = '  +1.11%  '

# Row 11
$ws.Range("D11").Value = '7.14'
$ws.Range("E11").Value = '  +2.65%  '

# Row 12
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  -0.16%  '

# Row 13
$ws.Range("D13").Value = '4.125.04'
$ws.Range("E13").Value = '  -1.21%  '

# Row 14
$ws.Range("E14").Value = '  +1.77%  '

# Row 15
$ws.Range("E15").Value = '  -0.53%  '

# Row 16
$ws.Range("E16").Value = '  +0.62%  '

# Row 17
$ws.Range("D17").Value = '3.526.38'
$ws.Range("E17").Value = '  -1.14%  '

# Row 18
$ws.Range("D18").Value = '65.117.21'
$ws.Range("E18").Value = '  +0.00%  '

# Row 19
$ws.Range("D19").Value = '10.11'
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("D20").Value = '14.37'
$ws.Range("E20").Value = '  -0.33%  '

# Row 21
$ws.Range("E21").Value = '  -3.11%  '

# Row 22
$ws.Range("D22").Value = '392.70'
$ws.Range("E22").Value = '  +1.06%  '

# Row 23
$ws.Range("D23").Value = '0.579'
$ws.Range("E23").Value = '  -0.17%  '

# Row 24
$ws.Range("D24").Value = '3.669.94'
$ws.Range("E24").Value = '  -1.18%  '

# Row 25
$ws.Range("D25").Value = '74.63'
$ws.Range("E25").Value = '  +0.55%  '

# Row 26
$ws.Range("E26").Value = '  +0.02%  '

# Row 27
$ws.Range("D27").Value = '0.0000111'
$ws.Range("E27").Value = '  -4.38%  '

# Row 28
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.60'
$ws.Range("E28").Value = '  +9.14%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.68'
$ws.Range("E29").Value = '  -0.46%  '

# Row 30
$ws.Range("E30").Value = '  -0.28%  '

# Row 32
$ws.Range("D32").Value = '8.33'

# Row 33
$ws.Range("D33").Value = '3.535.45'
$ws.Range("E33").Value = '  -1.18%  '

# Row 34
$ws.Range("E34").Value = '  +0.49%  '

# Row 35
$ws.Range("E35").Value = '  +0.01%  '

# Row 36
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("D37").Value = '5.30'
$ws.Range("E37").Value = '  +5.27%  '

# Row 38
$ws.Range("E38").Value = '  +1.67%  '

# Row 39
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("E40").Value = '  -0.48%  '

# Row 41
$ws.Range("D41").Value = '0.0806'
$ws.Range("E41").Value = '  -0.25%  '

# Row 42
$ws.Range("D42").Value = '0.822'
$ws.Range("E42").Value = '  -0.69%  '

# Row 43
$ws.Range("E43").Value = '  +4.66%  '

# Row 44
$ws.Range("D44").Value = '42.96'
$ws.Range("E44").Value = '  +0.67%  '

# Row 45
$ws.Range("D45").Value = '25.69'
$ws.Range("E45").Value = '  -5.23%  '

# Row 46
$ws.Range("E46").Value = '  +0.12%  '

# Row 47
$ws.Range("E47").Value = '  -1.09%  '

# Row 48
$ws.Range("E48").Value = '  +0.56%  '

# Row 49
$ws.Range("E49").Value = '  -0.34%  '

# Row 50
$ws.Range("D50").Value = '2.419.45'
$ws.Range("E50").Value = '  -3.19%  '

# Row 51
$ws.Range("E51").Value = '  +4.16%  '
